# New weekly price record for "Vega Modelo de Temuco - Haba".
# A new row is inserted above the current row 63, pushing the existing
# rows 63-66 down to 64-67, and the newly opened row 63 is filled with
# this week's data (same constant columns as the rest of the block, new
# date/volume/price figures).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 63:66 down to 64:67, opening up a blank row 63.
$ws.Rows("63:63").Insert()

# Populate the new row 63.
$ws.Cells.Item(63, 1).Value = 10
$ws.Cells.Item(63, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(63, 3).Value = "La Araucanía"
$ws.Cells.Item(63, 4).Value = 44585
$ws.Cells.Item(63, 5).Value = 9
$ws.Cells.Item(63, 6).Value = 100112026
$ws.Cells.Item(63, 7).Value = "Haba"
$ws.Cells.Item(63, 8).Value = "Sin especificar"
$ws.Cells.Item(63, 9).Value = "Primera"
$ws.Cells.Item(63, 10).Value = 55
$ws.Cells.Item(63, 11).Value = 20000
$ws.Cells.Item(63, 12).Value = 20000
$ws.Cells.Item(63, 13).Value = 20000
$ws.Cells.Item(63, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(63, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(63, 16).Value = 800
$ws.Cells.Item(63, 17).Value = 25
$ws.Cells.Item(63, 18).Value = "Hortaliza"
